# "Added cleaner stationary transforms (v0.14)"
#
# - Rename sheet "baseline-variables" -> "all-variables"
# - gdp row: source key (E2) A191RL1Q225SBEA -> GDPC1 (use the FRED real-GDP
#   series directly instead of the pre-computed SAAR % change series)
# - tdns1/tdns2/tdns3 (Treasury-FFR spread level/slope/curvature) stationary
#   transform (I23:I25) "d" (plain diff) -> "diff1" (the new, cleaner diff
#   transform)
# - view/selection bookkeeping: params sheet loses its scrolled-down
#   top-left cell, all-variables sheet's selection moves to E10, and
#   all-variables stays the active/visible tab.

$wb = $excel.ActiveWorkbook

$wsVars = $wb.Worksheets.Item("baseline-variables")

# Rename the sheet first so later lookups by new name also work if needed.
$wsVars.Name = "all-variables"

# --- data edits -----------------------------------------------------------

# gdp: sckey FRED series id
$wsVars.Range("E2").Value = "GDPC1"

# tdns1 / tdns2 / tdns3: stationary transform column (st)
$wsVars.Range("I23").Value = "diff1"
$wsVars.Range("I24").Value = "diff1"
$wsVars.Range("I25").Value = "diff1"

# --- view/selection edits --------------------------------------------------
# (params' scrolled-down topLeftCell is cleared automatically on save; no
# need to touch that sheet's view and risk stealing the active tab.)

# all-variables: move selection to E10 and leave it as the active tab
$wsVars.Activate() | Out-Null
$wsVars.Range("E10").Select() | Out-Null

Write-Output "done"
